# Scheduled GitHub Actions refresh of the cryptocurrency price/volume snapshot.
# Column D ("Price") holds values that can look numeric (e.g. "1.000", "0.9998");
# force text format before assigning so Excel does not silently coerce/round them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.795.39"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.59"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "277.35"
$ws.Range("E5").Value = "  -7.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5096"
$ws.Range("E7").Value = "  -4.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3516"
$ws.Range("E8").Value = "  -6.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.22"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06670"
$ws.Range("E10").Value = "  -7.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.99"
$ws.Range("E11").Value = "  -7.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8325"
$ws.Range("E12").Value = "  -6.08%  "
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.816.79"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.085"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.73"
$ws.Range("E16").Value = "  -5.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9996"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.11"
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008027"
$ws.Range("E19").Value = "  -5.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.831.18"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.725"
$ws.Range("E22").Value = "  -4.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.01"
$ws.Range("E23").Value = "  -6.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.080"
$ws.Range("E24").Value = "  -4.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.58"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.175"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.668"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.13"
$ws.Range("E28").Value = "  -5.12%  "
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.330"
$ws.Range("E30").Value = "  -8.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.242"
$ws.Range("E31").Value = "  -7.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08836"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04865"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7354"
$ws.Range("E34").Value = "  -7.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.885"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.150"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5236"
$ws.Range("E39").Value = "  -11.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.328"
$ws.Range("E40").Value = "  -11.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01846"
$ws.Range("E41").Value = "  -5.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9572"
$ws.Range("E42").Value = "  -10.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.12"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.202"
$ws.Range("E44").Value = "  -6.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.096"
$ws.Range("E45").Value = "  -9.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9993"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4591"
$ws.Range("E47").Value = "  -8.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1365"
$ws.Range("E48").Value = "  -8.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.49"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.198"
$ws.Range("E50").Value = "  -7.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.502"
$ws.Range("E51").Value = "  -7.22%  "
